$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A8").Value = (Get-Date -Year 2023 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B8").Value = "Internship"
$ws.Range("C8").Value = "Completed 8 hours assisting with daily operations"

$ws.Range("C9").Select()
